$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "B3" = 2
    "B11" = 2
    "B14" = 0
    "B16" = 0
    "B27" = 1
    "B28" = 2
    "B29" = 2
    "B32" = 2
    "B51" = 1
    "B60" = 1
    "B65" = 0
    "B66" = 0
    "B90" = 1
    "B92" = 2
    "B95" = 2
    "B98" = 1
    "B99" = 1
    "B105" = 1
    "B106" = 0
    "B117" = 2
    "B131" = 2
    "B136" = 2
    "B143" = 2
    "B150" = 1
    "B153" = 1
    "B154" = 0
    "B165" = 0
    "B167" = 1
    "B187" = 0
    "B191" = 1
    "B200" = 1
    "B204" = 2
    "B208" = 2
    "B209" = 1
    "B216" = 2
    "B228" = 1
    "B229" = 1
    "B236" = 2
    "B240" = 2
    "B243" = 2
    "B250" = 2
    "B253" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
